$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "64.328.95"
$ws.Range("E2").Value = "  -3.13%  "
# Row 3
Set-TextValue "D3" "3.176.55"
$ws.Range("E3").Value = "  -8.10%  "
# Row 4
$ws.Range("E4").Value = "  -0.02%  "
# Row 5
Set-TextValue "D5" "565.90"
$ws.Range("E5").Value = "  -3.78%  "
# Row 6
Set-TextValue "D6" "170.43"
$ws.Range("E6").Value = "  -3.38%  "
# Row 7
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
Set-TextValue "D7" "0.609"
$ws.Range("E7").Value = "  -0.59%  "
# Row 8
$ws.Range("B8").Value = "USDC"
$ws.Range("C8").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue "D8" "1.00"
$ws.Range("E8").Value = "  +0.03%  "
# Row 9
Set-TextValue "D9" "3.174.41"
$ws.Range("E9").Value = "  -8.14%  "
# Row 10
$ws.Range("E10").Value = "  -6.29%  "
# Row 11
$ws.Range("E11").Value = "  -4.56%  "
# Row 12
Set-TextValue "D12" "0.396"
$ws.Range("E12").Value = "  -4.88%  "
# Row 13
Set-TextValue "D13" "3.727.54"
$ws.Range("E13").Value = "  -8.11%  "
# Row 14
$ws.Range("E14").Value = "  +1.54%  "
# Row 15
Set-TextValue "D15" "27.43"
$ws.Range("E15").Value = "  -9.85%  "
# Row 16
Set-TextValue "D16" "64.327.42"
$ws.Range("E16").Value = "  -2.97%  "
# Row 17
Set-TextValue "D17" "0.0000164"
$ws.Range("E17").Value = "  -5.17%  "
# Row 18
Set-TextValue "D18" "3.175.77"
$ws.Range("E18").Value = "  -8.06%  "
# Row 19
$ws.Range("E19").Value = "  -3.69%  "
# Row 20
Set-TextValue "D20" "12.99"
$ws.Range("E20").Value = "  -5.94%  "
# Row 21
Set-TextValue "D21" "353.63"
$ws.Range("E21").Value = "  -5.44%  "
# Row 22
Set-TextValue "D22" "7.19"
$ws.Range("E22").Value = "  -5.69%  "
# Row 23
$ws.Range("E23").Value = "  +0.07%  "
# Row 24
Set-TextValue "D24" "69.06"
# Row 25
Set-TextValue "D25" "0.0000120"
$ws.Range("E25").Value = "  -4.45%  "
# Row 26
Set-TextValue "D26" "0.504"
# Row 27
Set-TextValue "D27" "9.56"
$ws.Range("E27").Value = "  -3.35%  "
# Row 28
$ws.Range("E28").Value = "  -1.05%  "
# Row 29
Set-TextValue "D29" "0.997"
$ws.Range("E29").Value = "  -0.27%  "
# Row 30
Set-TextValue "D30" "5.63"
$ws.Range("E30").Value = "  -4.50%  "
# Row 31
Set-TextValue "D31" "0.997"
$ws.Range("E31").Value = "  -0.23%  "
# Row 32
Set-TextValue "D32" "1.91"
$ws.Range("E32").Value = "  -4.92%  "
# Row 33
Set-TextValue "D33" "22.07"
$ws.Range("E33").Value = "  -7.01%  "
# Row 34
$ws.Range("E34").Value = "  -5.42%  "
# Row 35
Set-TextValue "D35" "6.65"
$ws.Range("E35").Value = "  -5.46%  "
# Row 36
$ws.Range("E36").Value = "  -6.78%  "
# Row 37
Set-TextValue "D37" "155.13"
$ws.Range("E37").Value = "  -3.32%  "
# Row 38
Set-TextValue "D38" "0.818"
$ws.Range("E38").Value = "  -7.60%  "
# Row 39
Set-TextValue "D39" "25.95"
$ws.Range("E39").Value = "  -8.80%  "
# Row 40
Set-TextValue "D40" "2.56"
$ws.Range("E40").Value = "  -1.94%  "
# Row 41
Set-TextValue "D41" "1.70"
$ws.Range("E41").Value = "  -6.27%  "
# Row 42
Set-TextValue "D42" "2.621.75"
$ws.Range("E42").Value = "  -5.04%  "
# Row 43
Set-TextValue "D43" "4.18"
$ws.Range("E43").Value = "  -7.20%  "
# Row 44
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D44" "6.04"
$ws.Range("E44").Value = "  -5.87%  "
# Row 45
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D45" "39.62"
$ws.Range("E45").Value = "  -1.05%  "
# Row 46
Set-TextValue "D46" "0.0658"
$ws.Range("E46").Value = "  -5.25%  "
# Row 47
Set-TextValue "D47" "23.88"
$ws.Range("E47").Value = "  -5.73%  "
# Row 48
Set-TextValue "D48" "323.43"
$ws.Range("E48").Value = "  -4.62%  "
# Row 49
$ws.Range("E49").Value = "  -7.43%  "
# Row 50
$ws.Range("E50").Value = "  -0.88%  "
# Row 51
Set-TextValue "D51" "1.00"
$ws.Range("E51").Value = "  -0.02%  "
